$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "Concepts" sheet as the very first tab.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$ws.Name = "Concepts"

# Header row
$ws.Range("A1").Value = "Concept ID"
$ws.Range("B1").Value = "Concept Name"
$ws.Range("C1").Value = "Definitions"
$ws.Range("D1").Value = "Example Usages"
$ws.Range("E1").Value = "Groups"
$ws.Range("F1").Value = "Synonyms"
$ws.Range("G1").Value = "Related Concepts"
$ws.Range("H1").Value = "References"

# Data rows
$ws.Range("A2").Value = "CON_0001"
$ws.Range("B2").Value = "Response Variable"
$ws.Range("F2").Value = "Dependent Variable;"

$ws.Range("A3").Value = "CON_0002"
$ws.Range("B3").Value = "Treatment"
$ws.Range("F3").Value = "Group;"

$ws.Range("A4").Value = "CON_0003"
$ws.Range("B4").Value = "Object"
$ws.Range("F4").Value = "Target Document;"

# The "Concept ID" column uses the same red-highlighted font the other
# id-style columns elsewhere in this workbook use.
$ws.Range("A1:A4").Font.Color = 255

# Column widths to roughly match the authored layout.
$ws.Columns.Item(1).ColumnWidth = 9.5 - (5 / 6)
$ws.Columns.Item(2).ColumnWidth = 15 - (5 / 6)
$ws.Columns.Item(3).ColumnWidth = 51.5 - (5 / 6)
$ws.Columns.Item(4).ColumnWidth = 51.5 - (5 / 6)
$ws.Columns.Item(5).ColumnWidth = 6.6666666666667 - (5 / 6)
$ws.Columns.Item(6).ColumnWidth = 16.6666666666667 - (5 / 6)
$ws.Columns.Item(7).ColumnWidth = 14.3333333333334 - (5 / 6)
$ws.Columns.Item(8).ColumnWidth = 15 - (5 / 6)

# Make it the active/selected sheet, matching the saved view state.
$ws.Activate()
$ws.Range("E9").Select()

# ---------------------------------------------------------------------------
# 2. Data-import fix on "Ontology Mapping": the "additionalProperty" column
#    (J) used the placeholder "no" -- flip every occurrence to ";" so the
#    new concept-driven class/property generation recognises it as an
#    (empty) multi-value list rather than a dead flag.
# ---------------------------------------------------------------------------
$om = $wb.Worksheets.Item("Ontology Mapping")
$lastRow = $om.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $om.Cells.Item($r, 10)
    if ($cell.Value() -eq "no") {
        $cell.Value = ";"
    }
}
